$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A8").Value = 2
$ws.Range("B8").Value = "avengers"

$ws.Range("B8").Select()
